$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Proyectos no se marca análisis por defecto", $true, $false, $false, $false, $false, $true, 1, $false, "Titles", 2)
